# Refactor the "Triggers" worksheet as "Processors".
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Triggers")
$ws.Name = "Processors"
